$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. Insert a new worksheet "2022-Q4" right after "总计", before "2022-Q3"
# ------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item(1)
$q4Sheet = $wb.Worksheets.Add($null, $totalSheet)
$q4Sheet.Name = "2022-Q4"
# Match the outline properties used by the sibling quarter sheets.
$q4Sheet.Outline.SummaryRow = 1
$q4Sheet.Outline.SummaryColumn = 1

# ------------------------------------------------------------------
# 2. Update the "总计" (summary) sheet: insert a new row for 2022-Q4
#    above the existing 2022-Q3 row, shifting everything else down.
# ------------------------------------------------------------------
$ws = $totalSheet

# Copy header-row style from A2 (existing data-row style) for the new A2.
$ws.Range("A5").Value = 3
$ws.Range("B5").Value = "2021-Q4"
$ws.Range("C5").Value = 1
$ws.Range("D5").Value = 0.01

$ws.Range("A4").Value = 2
$ws.Range("B4").Value = "2022-Q2"
$ws.Range("C4").Value = 10
$ws.Range("D4").Value = 0.8

$ws.Range("A3").Value = 1
$ws.Range("B3").Value = "2022-Q3"
$ws.Range("C3").Value = 6
$ws.Range("D3").Value = 0.35

$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "2022-Q4"
$ws.Range("C2").Value = 1
$ws.Range("D2").Value = 0.03

# Copy the style ("s=2") of column A down to the new A5 cell.
$ws.Range("A2").Copy()
$ws.Range("A5").PasteSpecial(-4122)

# ------------------------------------------------------------------
# 3. Populate the new "2022-Q4" sheet with the fund-holdings data.
# ------------------------------------------------------------------
$ws4 = $q4Sheet

$ws4.Range("B1").Value = "基金代码"
$ws4.Range("C1").Value = "基金名称"
$ws4.Range("D1").Value = "基金规模"
$ws4.Range("E1").Value = "股票总仓位"
$ws4.Range("F1").Value = "仓位占比"
$ws4.Range("G1").Value = "持有市值(亿元)"
$ws4.Range("H1").Value = "仓位排名"

$ws4.Range("A2").Value = 0

# B2:G2 hold text-typed values (fund code / name / percentages kept as
# literal text, matching the source data) -- force text storage via a
# "@" number format while writing, then drop back to the default style
# so the cells end up unstyled, same as the sibling quarter sheets.
$ws4.Range("B2:G2").NumberFormat = "@"
$ws4.Range("B2").Value = "200001"
$ws4.Range("C2").Value = "长城久恒灵活配置混合"
$ws4.Range("D2").Value = "0.85"
$ws4.Range("E2").Value = "94.35"
$ws4.Range("F2").Value = "3.70"
$ws4.Range("G2").Value = "0.0314"
$ws4.Range("B2:G2").Style = "Normal"

$ws4.Range("H2").Value = 1

# Copy the header style ("s=2") from the "总计" sheet's header row onto
# the new sheet's header cells and the A2 index cell.
$totalSheet.Range("B1").Copy()
$ws4.Range("B1:H1").PasteSpecial(-4122)

$totalSheet.Range("A2").Copy()
$ws4.Range("A2").PasteSpecial(-4122)

# Restore the original active sheet / selection state ("总计" was active).
$totalSheet.Activate() | Out-Null
$totalSheet.Range("A1").Select() | Out-Null
